$wb = $excel.ActiveWorkbook
$before = $wb.Worksheets.Item("HokenMst")
$newws = $wb.Worksheets.Add($before)
$newws.Name = "HpInf"

# Header row (left to right)
$newws.Cells.Item(1,1).Value = "hp_id"
$newws.Cells.Item(1,2).Value = "start_date"
$newws.Cells.Item(1,3).Value = "hp_cd"
$newws.Cells.Item(1,4).Value = "rousai_hp_cd"
$newws.Cells.Item(1,5).Value = "hp_name"
$newws.Cells.Item(1,6).Value = "rece_hp_name"
$newws.Cells.Item(1,7).Value = "kaisetu_name"
$newws.Cells.Item(1,8).Value = "post_cd"
$newws.Cells.Item(1,9).Value = "pref_no"
$newws.Cells.Item(1,10).Value = "address1"
$newws.Cells.Item(1,11).Value = "address2"
$newws.Cells.Item(1,12).Value = "tel"
$newws.Cells.Item(1,13).Value = "create_date"
$newws.Cells.Item(1,14).Value = "create_id"
$newws.Cells.Item(1,15).Value = "create_machine"
$newws.Cells.Item(1,16).Value = "update_date"
$newws.Cells.Item(1,17).Value = "update_id"
$newws.Cells.Item(1,18).Value = "update_machine"
$newws.Cells.Item(1,19).Value = "fax_no"
$newws.Cells.Item(1,20).Value = "other_contacts"

# Data row 2
$newws.Cells.Item(2,1).Value = 998
$newws.Cells.Item(2,2).Value = 0
$newws.Cells.Item(2,4).Value = "sfdffsj"
$newws.Cells.Item(2,5).Value = "sfdffsjfklsjrpoiqewrejksdfjalkjfdjfqwoiejfljlskdjfldsjflsdjfljfoiwjeoijoiewjiore"
$newws.Cells.Item(2,6).Value = "sfdffsjfklsjrpoiqewrejksdfjalkjfdjfqwoiejfljlskdjfldsjflsdjfljfoiwjeoijoiewjiore"
$newws.Cells.Item(2,7).Value = "sfdffsjfklsjrpoiqewrejksdfjalkjfdjfqwoie"
$newws.Cells.Item(2,8).Value = 1231232
$newws.Cells.Item(2,9).Value = 17
$newws.Cells.Item(2,10).Value = "sfdffsjfklsjrpoiqewrejksdfjalkjfdjfqwoiejfljlskdjfldsjflsdjfljfoiwjeoijoiewjiorewjroiewjroiewjroiewj"
$newws.Cells.Item(2,11).Value = "sfdffsjfklsjrpoiqewrejksdfjalkjfdjfqwoiejfljlskdjfldsjflsdjfljfoiwjeoijoiewjiorewjroiewjroiewjroiewj"
$newws.Cells.Item(2,12).Value = "sfdffsjfklsjrpo"
$newws.Cells.Item(2,3).Value = "abcd"
$newws.Cells.Item(2,13).Value = 45040.84794363426
$newws.Cells.Item(2,13).NumberFormat = "mm:ss.0"
$newws.Cells.Item(2,14).Value = 2
$newws.Cells.Item(2,16).Value = 45184.049576296296
$newws.Cells.Item(2,16).NumberFormat = "mm:ss.0"
$newws.Cells.Item(2,17).Value = 2

$hoken = $wb.Worksheets.Item("HokenMst")
$hoken.Range("B2").Select()

$newws.Range("J12").Select()
